$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings. Some new values would otherwise be auto-
# coerced to numbers by Excel (dropping meaningful trailing zeros, e.g.
# "11.00" -> 11), so force text format first for those cells only.
$ws.Range("D2").Value = "69.770.07"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "3.693.37"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "674.15"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.38"
$ws.Range("E6").Value = "  +2.48%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.11"
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.444"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000234"
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.32"
$ws.Range("E13").Value = "  +3.60%  "
$ws.Range("D14").Value = "3.692.39"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").Value = "69.717.32"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "16.20"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.52"
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "472.48"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.81"
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.650"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "80.15"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").Value = "3.841.89"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("E24").Value = "  +7.48%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.00"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.14"
$ws.Range("E27").Value = "  +2.00%  "
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.03"
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.168"
$ws.Range("E31").Value = "  +5.83%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.54"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.90"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("D35").Value = "3.689.60"
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("E36").Value = "  +4.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.15"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.29"
$ws.Range("E39").Value = "  +4.31%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0911"
$ws.Range("E41").Value = "  +1.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "175.87"
$ws.Range("E42").Value = "  +2.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.936"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "46.98"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("E45").Value = "  +4.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "27.98"
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.08"
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("E50").Value = "  +2.18%  "
$ws.Range("E51").Value = "  +0.36%  "
